# Logged Week 15 and simulated Week 16
# Appends this week's per-play/per-kick data to the running season logs
# (stored as space-separated numbers inside shared-string cells) and bumps
# the season-to-date totals on the summary sheets accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: running logs of per-play yardage/result values.
# Each cell holds a single space-separated string of numbers that grows
# by one week's worth of entries; we must append, not overwrite.
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

function Append-Series([object]$ws, [string]$cellRef, [int[]]$values) {
    $current = $ws.Range($cellRef).Value2
    $addition = ($values -join " ")
    if ([string]::IsNullOrEmpty($current)) {
        $ws.Range($cellRef).Value2 = $addition
    } else {
        $ws.Range($cellRef).Value2 = "$current $addition"
    }
}

Append-Series $wsYDS "B2" @(3,6,2,3,5,2,-1,6,6,3,8,0,3,12,30,2,3,2,6,0,17)
Append-Series $wsYDS "B3" @(5,7,7,-2,3,23,5,8,1,5,5,6,14,19,7,5,8,5,16,6,10,3,6,18,11,13)
Append-Series $wsYDS "C2" @(8,-1,4,-3,5,-2,1,1,9,2,2,6,6,3,-1,3,8,0,0,7,3,-1,4,1,1,-1,3,-3,1,0,-5)
Append-Series $wsYDS "C3" @(8,9,7,40,6,7,33,8,4,3,17,5,7)

# ---------------------------------------------------------------------
# OFF sheet: season-to-date offensive totals.
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 163
$wsOFF.Range("F2").Value = 19
$wsOFF.Range("G2").Value = 39
$wsOFF.Range("I2").Value = 2
$wsOFF.Range("J2").Value = 28
$wsOFF.Range("N2").Value = 13
$wsOFF.Range("O2").Value = 14

$wsOFF.Range("C3").Value = 187
$wsOFF.Range("E3").Value = 21
$wsOFF.Range("F3").Value = 142
$wsOFF.Range("G3").Value = 53
$wsOFF.Range("I3").Value = 64
$wsOFF.Range("J3").Value = 56
$wsOFF.Range("L3").Value = 350
$wsOFF.Range("M3").Value = 239
$wsOFF.Range("Q3").Value = 525

# ---------------------------------------------------------------------
# DEF sheet: season-to-date defensive totals.
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 140
$wsDEF.Range("E2").Value = 7
$wsDEF.Range("F2").Value = 44
$wsDEF.Range("G2").Value = 37
$wsDEF.Range("H2").Value = 7
$wsDEF.Range("I2").Value = 6
$wsDEF.Range("J2").Value = 20
$wsDEF.Range("N2").Value = 23

$wsDEF.Range("B3").Value = 5
$wsDEF.Range("C3").Value = 199
$wsDEF.Range("E3").Value = 29
$wsDEF.Range("F3").Value = 110
$wsDEF.Range("G3").Value = 41
$wsDEF.Range("H3").Value = 21
$wsDEF.Range("I3").Value = 52
$wsDEF.Range("J3").Value = 61
$wsDEF.Range("L3").Value = 290
$wsDEF.Range("M3").Value = 191
$wsDEF.Range("Q3").Value = 505

# ---------------------------------------------------------------------
# ST sheet: season-to-date special-teams totals plus running logs of
# per-kick distance/return numbers (same append pattern as YDS).
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 81
$wsST.Range("D2").Value = 51
$wsST.Range("L2").Value = 116
$wsST.Range("B3").Value = 67

Append-Series $wsST "D3" @(40,53,48,45,39,42,30,39)
Append-Series $wsST "D4" @(0,34,0,0,4,8,0,0)
Append-Series $wsST "D5" @(6,11,3,2,9,10,0,0,0)
Append-Series $wsST "B6" @(17,22)

# ---------------------------------------------------------------------
# TURNS sheet: season-to-date turnover totals.
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B3").Value = 8
$wsTURNS.Range("D3").Value = 7
$wsTURNS.Range("E3").Value = 9

# ---------------------------------------------------------------------
# PEN sheet: season-to-date penalty totals.
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B2").Value = 15
